$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.681.69'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '2.477.20'
$ws.Range('E3').Value = '  +0.70%  '
$ws.Range('E4').Value = '  +0.07%  '
$cell = $ws.Range('D5')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '319.47'
$cell.Style = $origStyle
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('E7').Value = '  +0.93%  '
$ws.Range('E8').Value = '  +0.04%  '
$cell = $ws.Range('D9')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.512'
$cell.Style = $origStyle
$ws.Range('E9').Value = '  +0.47%  '
$cell = $ws.Range('D10')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0862'
$cell.Style = $origStyle
$ws.Range('E10').Value = '  +8.11%  '
$cell = $ws.Range('D11')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '33.07'
$cell.Style = $origStyle
$ws.Range('E11').Value = '  +2.22%  '
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('D13').Value = '2.859.65'
$ws.Range('E13').Value = '  +0.80%  '
$cell = $ws.Range('D14')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '6.89'
$cell.Style = $origStyle
$ws.Range('E14').Value = '  +0.78%  '
$cell = $ws.Range('D15')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '15.52'
$cell.Style = $origStyle
$ws.Range('E15').Value = '  -1.63%  '
$ws.Range('D16').Value = '2.462.46'
$ws.Range('E16').Value = '  +0.27%  '
$cell = $ws.Range('D17')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.794'
$cell.Style = $origStyle
$ws.Range('E17').Value = '  +2.48%  '
$ws.Range('D18').Value = '41.605.84'
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('E19').Value = '  -0.28%  '
$ws.Range('D20').Value = '0.0₃0942'
$ws.Range('E20').Value = '  +0.79%  '
$cell = $ws.Range('D21')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '70.63'
$cell.Style = $origStyle
$ws.Range('E21').Value = '  -0.07%  '
$cell = $ws.Range('D22')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '11.25'
$cell.Style = $origStyle
$ws.Range('E22').Value = '  -0.37%  '
$cell = $ws.Range('D23')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '239.92'
$cell.Style = $origStyle
$ws.Range('E23').Value = '  +0.72%  '
$ws.Range('E24').Value = '  +1.73%  '
$ws.Range('E25').Value = '  +2.37%  '
$cell = $ws.Range('D27')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '25.00'
$cell.Style = $origStyle
$ws.Range('E27').Value = '  +2.99%  '
$cell = $ws.Range('D28')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.25'
$cell.Style = $origStyle
$ws.Range('E28').Value = '  -0.33%  '
$cell = $ws.Range('D29')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '9.76'
$cell.Style = $origStyle
$ws.Range('E29').Value = '  +0.62%  '
$cell = $ws.Range('D30')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '36.70'
$cell.Style = $origStyle
$ws.Range('E30').Value = '  +4.67%  '
$cell = $ws.Range('D31')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '156.98'
$cell.Style = $origStyle
$ws.Range('E31').Value = '  +0.65%  '
$cell = $ws.Range('D32')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '5.43'
$cell.Style = $origStyle
$ws.Range('E32').Value = '  -0.33%  '
$ws.Range('E33').Value = '  -0.03%  '
$cell = $ws.Range('D34')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0764'
$cell.Style = $origStyle
$ws.Range('E34').Value = '  +0.62%  '
$ws.Range('E35').Value = '  -0.25%  '
$ws.Range('E36').Value = '  -1.61%  '
$cell = $ws.Range('D37')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.117'
$cell.Style = $origStyle
$ws.Range('E37').Value = '  +1.89%  '
$ws.Range('E38').Value = '  +2.84%  '
$cell = $ws.Range('D39')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.89'
$cell.Style = $origStyle
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('E40').Value = '  +1.37%  '
$cell = $ws.Range('D41')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '4.01'
$cell.Style = $origStyle
$ws.Range('E41').Value = '  +1.81%  '
$ws.Range('E42').Value = '  -1.53%  '
$ws.Range('D43').Value = '2.001.61'
$ws.Range('E43').Value = '  +1.55%  '
$ws.Range('E44').Value = '  +0.79%  '
$cell = $ws.Range('D45')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '18.72'
$cell.Style = $origStyle
$ws.Range('E45').Value = '  -1.44%  '
$cell = $ws.Range('D46')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.98'
$cell.Style = $origStyle
$ws.Range('E46').Value = '  +2.33%  '
$cell = $ws.Range('D47')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '9.54'
$cell.Style = $origStyle
$ws.Range('E47').Value = '  +5.76%  '
$ws.Range('D48').Value = '2.715.55'
$ws.Range('E48').Value = '  +0.79%  '
$cell = $ws.Range('D49')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '97.70'
$cell.Style = $origStyle
$ws.Range('E49').Value = '  +0.88%  '
$cell = $ws.Range('D50')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '75.63'
$cell.Style = $origStyle
$ws.Range('E50').Value = '  +5.27%  '
$cell = $ws.Range('D51')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '67.35'
$cell.Style = $origStyle
$ws.Range('E51').Value = '  +0.97%  '
